$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position bookkeeping (best-effort; mirrors the yWindow move in workbook.xml) ---
try { $excel.ActiveWindow.Top = 912 } catch { }

# --- Header row (row 1): replace the old title with column headers, all bold ---
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Population"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Source"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.Underline = $false

# --- Data row (row 2): new numeric + text data ---
$ws.Range("A2").Value = 2500
$ws.Range("B2").Value = 90000000000
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("C2").Value = "Everywhere"
$ws.Range("D2").Value = "Nothing at all"

# --- Old header row (row 3) is no longer needed: clear its contents/format ---
$ws.Range("A3:D3").ClearContents()
$ws.Range("A3:D3").Style = "Normal"

# --- Move the active selection ---
$ws.Range("H9").Select()

Write-Output "done"
